$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row (D=Fecha, K=Variedad, L=Calidad, M=Volumen, N=Precio minimo,
# O=Precio maximo, P=Precio promedio ponderado, R=Origen, S=Precio $/Kg)

$ws.Range("D2").Value = 44355
$ws.Range("K2").Value = "Mankaki"
$ws.Range("M2").Value = 270
$ws.Range("R2").Value = "Región Metropolitana"

$ws.Range("D3").Value = 44342
$ws.Range("L3").Value = "Primera"

$ws.Range("D4").Value = 44313
$ws.Range("M4").Value = 270
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("S4").Value = 1194

$ws.Range("D5").Value = 44305
$ws.Range("M5").Value = 250
$ws.Range("N5").Value = 24000
$ws.Range("O5").Value = 25000
$ws.Range("P5").Value = 24500
$ws.Range("R5").Value = "Región de O'Higgins"
$ws.Range("S5").Value = 1361

$ws.Range("D6").Value = 44301
$ws.Range("K6").Value = "Hachiya"
$ws.Range("L6").Value = "Segunda"
$ws.Range("M6").Value = 250
$ws.Range("N6").Value = 20000
$ws.Range("O6").Value = 21000
$ws.Range("P6").Value = 20500
$ws.Range("S6").Value = 1139
